# Update SanityCheck for DataModeler
# Removes the leftover "title" row (row 1) from each of the three sheets,
# renames two relation-description strings on the SalesOrderItem sheet,
# and refreshes the bestFit column widths / cell selections that Excel
# recalculates as a consequence of those edits.

$wb = $excel.ActiveWorkbook

# Helper: ColumnWidth must be set to a value that this engine's
# chars-to-stored-width rounding (nearest 1/6 column-width unit) maps
# back to the desired stored width as closely as possible.
function Set-BestFitColumnWidth($col, $targetStoredWidth) {
    $n = [Math]::Round($targetStoredWidth * 6 - 5)
    $cw = $n / 6.0
    $col.ColumnWidth = $cw
}

# ---- Sheet "SalesOrder" --------------------------------------------------
$ws1 = $wb.Worksheets.Item("SalesOrder")
$ws1.Range("A1").ClearContents()

Set-BestFitColumnWidth $ws1.Columns.Item(1) 7.42578125
Set-BestFitColumnWidth $ws1.Columns.Item(4) 8.140625
Set-BestFitColumnWidth $ws1.Columns.Item(5) 8.85546875

# ---- Sheet "SalesOrderItem" ----------------------------------------------
$ws2 = $wb.Worksheets.Item("SalesOrderItem")
$ws2.Range("A1").ClearContents()
$ws2.Range("C2").Value = "SalesOrder.SalesOrderItem.SalesOrderItem"
$ws2.Range("D2").Value = "SalesOrderItem.RelationName.Product"

Set-BestFitColumnWidth $ws2.Columns.Item(1) 4.28515625
Set-BestFitColumnWidth $ws2.Columns.Item(2) 8.7109375
Set-BestFitColumnWidth $ws2.Columns.Item(3) 40.5703125
Set-BestFitColumnWidth $ws2.Columns.Item(4) 31

# ---- Sheet "Product" -------------------------------------------------
$ws3 = $wb.Worksheets.Item("Product")
$ws3.Range("A1").ClearContents()

# Reset sheet3's stale selection back to the default (A1) - must happen
# before we restore focus to SalesOrderItem below, since selecting a
# range also activates that sheet.
$ws3.Range("A1").Select()

# ---- Final selections / active sheet -------------------------------------
$ws1.Range("C10").Select()

# SalesOrderItem is the workbook's active tab; select it last so it stays
# the active sheet, and leave the cursor on D3 as in the edited file.
$ws2.Range("D3").Select()
